$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Insert 3 new rows before the current totals row (row 13), pushing it down to row 16
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()

# Update Status of rows 7-9 from "Aperto" to "Chiuso"
$ws.Cells.Item(7,4).Value2 = "Chiuso"
$ws.Cells.Item(8,4).Value2 = "Chiuso"
$ws.Cells.Item(9,4).Value2 = "Chiuso"

# Update Owner of row 9 from "?" to "Marco"
$ws.Cells.Item(9,6).Value2 = "Marco"

# Fill in new row 13: Clessidra / Logica / 2 / Chiuso / No / Marco
$ws.Cells.Item(13,1).Value2 = "Clessidra"
$ws.Cells.Item(13,2).Value2 = "Logica"
$ws.Cells.Item(13,3).Value2 = 2
$ws.Cells.Item(13,4).Value2 = "Chiuso"
$ws.Cells.Item(13,5).Value2 = "No"
$ws.Cells.Item(13,6).Value2 = "Marco"

# Fill in new row 14: Clessidra / Assets / 1 / Chiuso / No / Daniele
$ws.Cells.Item(14,1).Value2 = "Clessidra"
$ws.Cells.Item(14,2).Value2 = "Assets"
$ws.Cells.Item(14,3).Value2 = 1
$ws.Cells.Item(14,4).Value2 = "Chiuso"
$ws.Cells.Item(14,5).Value2 = "No"
$ws.Cells.Item(14,6).Value2 = "Daniele"

# Fill in new row 15: Spawn tagli / Debug / 2 / Aperto / Sì / Marco
$ws.Cells.Item(15,1).Value2 = "Spawn tagli"
$ws.Cells.Item(15,2).Value2 = "Debug"
$ws.Cells.Item(15,3).Value2 = 2
$ws.Cells.Item(15,4).Value2 = "Aperto"
$ws.Cells.Item(15,5).Value2 = "Sì"
$ws.Cells.Item(15,6).Value2 = "Marco"

# Resize the table to include the new rows (A1:F16), totals row becomes row 16
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F16"))

# Update the totals row formula for the Effort column (custom calculation)
$ws.Cells.Item(16,3).Formula = '=SUM([Effort]) - SUMIF([Status],$D$7,[Effort])'

# Update the selected cell
$ws.Range("D1").Select()
